$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@("DtqQ8kgAAAAJ", "Comparative physiology")
    ,@("gAYL31cAAAAJ", "Developmental biology | physiology | neuroscience")
    ,@("kBqFcdAAAAAJ", $null)
    ,@("-DHZ-BMAAAAJ", "Evolution | Ecology | Genetics/Genomics | Fish")
    ,@("4T2FMrgAAAAJ", "animal welfare | social behavior | animal personality | precision livestock farming")
    ,@("0sXoS44AAAAJ", "toxicology | physiology | pharmacology | molecular biology | metabolism")
    ,@("O17cHX0AAAAJ", $null)
    ,@("9I87lFkAAAAJ", "Reproductive Physiology")
    ,@("T-2je90AAAAJ", "Marine Mammal Ecology")
    ,@("mZ68FaQAAAAJ", "Parasitology | Arthropods | Genomics")
    ,@("_3VcPE0AAAAJ", "physiology | evolution | carryover effects | fetal programming")
    ,@("UGaV_2QAAAAJ", "Musculoskeletal disease | Preclinical Imaging | Objective measurement of animal pain | drug development")
    ,@("jfGWJdAAAAAJ", "biology | neuroscience | stress | life course health")
    ,@("ScMq2UkAAAAJ", "Biology | Ecological Physiology | Comparative Physiology | Migration | Animal Energetics")
    ,@("UTqnX-sAAAAJ", "Vertebrate paleontology | Evolutionary developmental biology")
    ,@("-Udc0ekAAAAJ", "Evolutionary Physiology | Behavioural ecology | Quantitative Genetics | Energetics | Personality")
    ,@("iPFpoPIAAAAJ", "Arctic | Ecology | Ecotoxicology | Environmental chemistry | Wildlife")
    ,@("TzNE64UAAAAJ", "stress | ecophysiology | hormones | neuroscience | microbiome")
    ,@("yHv0N68AAAAJ", "insect molecular physiology and neuroendocrinology")
    ,@("xhf9FkoAAAAJ", "genomics | behavior | social insects | neuroscience | evolution")
    ,@("Y11pVRYAAAAJ", "Sexual selection | evolutionary biology | animal behaviour | behavioural ecology")
    ,@("ErzXye8AAAAJ", "Marine evolutionary ecology | Molecular ecology")
    ,@("gXBATE4AAAAJ", "ecophysiology | foraging | predator avoidance | social behaviour | fishing-induced evolution")
    ,@("ekF5Q_UAAAAJ", "Animal welfare | Animal behaviour")
    ,@("QL99sHsAAAAJ", "learning health systems | artificial intelligence | multimorbidity | patient reported outcome | ehealth")
    ,@("IXlyPYkAAAAJ", "environmental toxicology | exposure science | ecotoxicology | environmental health | mercury")
    ,@("w2cncq4AAAAJ", "fish health and disease | aquaculture | biosecurity")
    ,@("3RkQotcAAAAJ", "Toxicology | Animal Physiology | Genomics")
    ,@("GzeHDMwAAAAJ", "ecotoxicology | birds | pesticides | agroecosystems | rivers and wetlands")
    ,@("ei5sM98AAAAJ", "social organization | insect communication | evolution of cooperation | honey bees")
    ,@("l3u-vkoAAAAJ", "Evolutionary physiology")
    ,@("_rT27wwAAAAJ", "Reproductive endocrinology | metabolism | transcription")
    ,@("od8jUdcAAAAJ", $null)
    ,@("XCWmQ6gAAAAJ", "Evolutionary Ecology | Reintroduction Biology | Species at Risk | Reproductive Ecology")
    ,@("o8U1CKwAAAAJ", "Virology | oncolytic immunotherapy | viral vectored vaccines | gene therapy")
    ,@("fjeFgDwAAAAJ", "aquaporin | ion transport | ammonia | insect | septate junctions")
    ,@("UawkFg4AAAAJ", "fish health | environmental impacts | environmental monitoring | sculpin | water quality")
    ,@("6yGqzVoAAAAJ", "Stem cells and regeneration | muscle satellite cells | myogenesis and adipogenesis | brown adipocytes")
    ,@("AWK2FDoAAAAJ", "energetics | ecology | behaviour | hibernation | bats")
    ,@("CAlipoMAAAAJ", "Fisheries | Fish Ecology | Conservation Science | Conservation Physiology | Freshwater Ecology")
    ,@("_uWapBgAAAAJ", "Neuroscience  Synaptic Transmission  Exocytosis  Addiction  Stxbp1")
    ,@("4YaeoLQAAAAJ", "Origins & Evolution of Animal Body Plan Diversity")
    ,@("6ncgZ2AAAAAJ", "Wildlife Biology | Traditional Food Security")
    ,@("FuVV0UYAAAAJ", "Biological psychiatry")
    ,@("OnCLeigAAAAJ", "Animal Behaviour | Animal Welfare | Applied Ethology")
    ,@("V3M92AQAAAAJ", "neuromuscular control")
    ,@("oc7XRMYAAAAJ", "vertebrate paleontology | evolution | developmental evolution | Archosauria")
    ,@("l2C58-wAAAAJ", "synaptic physiology | GPCRs | optogenetics")
    ,@("Ao3Hw_YAAAAJ", "conservation | physiology | policy")
    ,@("8P90kTYAAAAJ", "Behavioural Neuroendocrinology | Behavioural Ecology | Neuroscience")
    ,@("xGL-ao0AAAAJ", "Neurobiology")
    ,@("nph1nnIAAAAJ", "genetics | metabolism | C. elegans | membranes")
    ,@("LT35pKEAAAAJ", "developmental biology")
    ,@("LIcw3QgAAAAJ", "Neuroscience | neuromodulation | respiratory control | development")
    ,@("zaiIL4oAAAAJ", "ecoimmunology | ecological immunology | comparative stress physiology | parasitic manipulation of behaviour")
    ,@("JR_Dml4AAAAJ", "cloud physics medical physics")
    ,@("tqzOiyoAAAAJ", $null)
    ,@("ImUDvyMAAAAJ", "Vector Biology | Population Genetics/Genomics | Ecology")
    ,@("XwfzBI4AAAAJ", "Evolutionary Genetics")
    ,@("ICJSbLEAAAAJ", "Pigmentation | Geographic Variation | Carotenoids | Birds")
    ,@("esIT6IEAAAAJ", "Sexual Selection")
    ,@("OVOrQ4EAAAAJ", $null)
    ,@("wdDq5nkAAAAJ", "PI | Dog Aging Project")
    ,@("4OsbXUcAAAAJ", "Community Ecology | behavioral ecology | conservation ecology")
)

$startRow = 412
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $data[$i][0]
    if ($null -ne $data[$i][1]) {
        $ws.Cells.Item($r, 2).Value = $data[$i][1]
    }
}

Write-Output "Added $($data.Count) rows"